# Generate Report for handoff
# Two new files entered the localization pipeline:
#   - 6172aa8c-fd0c-45bd-b147-5e3b6aba0a9e.md  (now "Ready for handoff")
#   - 97776504-99c2-46cd-96fa-5e598604d67d.md  (now "Ready for handoff")
# The two files that used to be "Ready for handoff" are now "In Translation".
# The handoff timestamps on the overview rows were refreshed as well.

$wb = $excel.ActiveWorkbook

$e2eCommit = "5bccc3e767d7e81902e1ace4d53aff6088f406de"
$zhCommit  = "5aeb5130bd313757708da7d154cead4a8de28737"
$deCommit  = "d11f565a8e119a13b4fb0c5813977056833c96c2"

function E2eUrl($fn) {
  return "https://github.com/OpenLocalizationTest/oltest/blob/$e2eCommit/e2e/$fn"
}
function ConfigUrl() {
  return "https://github.com/OpenLocalizationTest/oltest/blob/$e2eCommit/.localization-config"
}
function ZhUrl($fn) {
  return "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$fn"
}
function DeUrl($fn) {
  return "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$fn"
}

$file98dc = "98dc3ec0-bd75-4257-a4a3-b6d9775752e6.md"
$fileF0b8 = "f0b871c6-8dda-4742-b2b2-fc540f9b6ff1.md"
$file6172 = "6172aa8c-fd0c-45bd-b147-5e3b6aba0a9e.md"
$file9777 = "97776504-99c2-46cd-96fa-5e598604d67d.md"
$fileConfig = ".localization-config"

$xlf98dcZh = "98dc3ec0-bd75-4257-a4a3-b6d9775752e6.9326827119ce91f2b6b7bfacbc27f0b28571abda.zh-cn.xlf"
$xlfF0b8Zh = "f0b871c6-8dda-4742-b2b2-fc540f9b6ff1.094af1df533f843959ba78c63ad5583f3e88a0b8.zh-cn.xlf"
$xlf6172Zh = "6172aa8c-fd0c-45bd-b147-5e3b6aba0a9e.bc38ad4e19d5ef4e2011c99ee54de0a89e77b126.zh-cn.xlf"
$xlf9777Zh = "97776504-99c2-46cd-96fa-5e598604d67d.edc3a0cfff62900734d2c7bc2c7b34f624764c99.zh-cn.xlf"

$xlf98dcDe = "98dc3ec0-bd75-4257-a4a3-b6d9775752e6.9326827119ce91f2b6b7bfacbc27f0b28571abda.de-de.xlf"
$xlfF0b8De = "f0b871c6-8dda-4742-b2b2-fc540f9b6ff1.094af1df533f843959ba78c63ad5583f3e88a0b8.de-de.xlf"
$xlf6172De = "6172aa8c-fd0c-45bd-b147-5e3b6aba0a9e.bc38ad4e19d5ef4e2011c99ee54de0a89e77b126.de-de.xlf"
$xlf9777De = "97776504-99c2-46cd-96fa-5e598604d67d.edc3a0cfff62900734d2c7bc2c7b34f624764c99.de-de.xlf"

$zhTimestamp = "2016-01-13 04:07:03"
$deTimestamp = "2016-01-13 04:07:25"
$epoch = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Cells.Item(2,2).Value = "In Translation"
$ws1.Cells.Item(2,3).Value = "In Translation"
$ws1.Cells.Item(3,2).Value = "In Translation"
$ws1.Cells.Item(3,3).Value = "In Translation"

$ws1.Cells.Item(4,1).Value = $file6172
$ws1.Cells.Item(4,2).Value = "Ready for handoff"
$ws1.Cells.Item(4,3).Value = "Ready for handoff"

$ws1.Cells.Item(5,1).Value = $file9777
$ws1.Cells.Item(5,2).Value = "Ready for handoff"
$ws1.Cells.Item(5,3).Value = "Ready for handoff"

$ws1.Cells.Item(6,1).Value = $fileConfig
$ws1.Cells.Item(6,2).Value = "Not to be localized"
$ws1.Cells.Item(6,3).Value = "Not to be localized"

$ws1.Cells.Item(1,1).Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Cells.Item(2,1), (E2eUrl $file98dc), "", "", $file98dc) | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(3,1), (E2eUrl $fileF0b8), "", "", $fileF0b8) | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(4,1), (E2eUrl $file6172), "", "", $file6172) | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(5,1), (E2eUrl $file9777), "", "", $file9777) | Out-Null
$ws1.Hyperlinks.Add($ws1.Cells.Item(6,1), (ConfigUrl), "", "", $fileConfig) | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Cells.Item(2,2).Value = "In Translation"
$ws2.Cells.Item(3,2).Value = "In Translation"
$ws2.Cells.Item(2,4).Value = $zhTimestamp
$ws2.Cells.Item(3,4).Value = $zhTimestamp

$ws2.Cells.Item(4,1).Value = $file6172
$ws2.Cells.Item(4,2).Value = "Ready for handoff"
$ws2.Cells.Item(4,3).Value = $xlf6172Zh
$ws2.Cells.Item(4,4).Value = $zhTimestamp
$ws2.Cells.Item(4,7).Value = $epoch
$ws2.Cells.Item(4,8).Value = "Include"

$ws2.Cells.Item(5,1).Value = $file9777
$ws2.Cells.Item(5,2).Value = "Ready for handoff"
$ws2.Cells.Item(5,3).Value = $xlf9777Zh
$ws2.Cells.Item(5,4).Value = $zhTimestamp
$ws2.Cells.Item(5,7).Value = $epoch
$ws2.Cells.Item(5,8).Value = "Include"

$ws2.Cells.Item(6,1).Value = $fileConfig
$ws2.Cells.Item(6,2).Value = "Not to be localized"
$ws2.Cells.Item(6,4).Value = $epoch
$ws2.Cells.Item(6,7).Value = $epoch
$ws2.Cells.Item(6,8).Value = "Ignored"

$ws2.Cells.Item(1,1).Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Cells.Item(2,1), (E2eUrl $file98dc), "", "", $file98dc) | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(2,3), (ZhUrl $xlf98dcZh), "", "", $xlf98dcZh) | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(3,1), (E2eUrl $fileF0b8), "", "", $fileF0b8) | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(3,3), (ZhUrl $xlfF0b8Zh), "", "", $xlfF0b8Zh) | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(4,1), (E2eUrl $file6172), "", "", $file6172) | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(4,3), (ZhUrl $xlf6172Zh), "", "", $xlf6172Zh) | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(5,1), (E2eUrl $file9777), "", "", $file9777) | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(5,3), (ZhUrl $xlf9777Zh), "", "", $xlf9777Zh) | Out-Null
$ws2.Hyperlinks.Add($ws2.Cells.Item(6,1), (ConfigUrl), "", "", $fileConfig) | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Cells.Item(2,2).Value = "In Translation"
$ws3.Cells.Item(3,2).Value = "In Translation"
$ws3.Cells.Item(2,4).Value = $deTimestamp
$ws3.Cells.Item(3,4).Value = $deTimestamp

$ws3.Cells.Item(4,1).Value = $file6172
$ws3.Cells.Item(4,2).Value = "Ready for handoff"
$ws3.Cells.Item(4,3).Value = $xlf6172De
$ws3.Cells.Item(4,4).Value = $deTimestamp
$ws3.Cells.Item(4,7).Value = $epoch
$ws3.Cells.Item(4,8).Value = "Include"

$ws3.Cells.Item(5,1).Value = $file9777
$ws3.Cells.Item(5,2).Value = "Ready for handoff"
$ws3.Cells.Item(5,3).Value = $xlf9777De
$ws3.Cells.Item(5,4).Value = $deTimestamp
$ws3.Cells.Item(5,7).Value = $epoch
$ws3.Cells.Item(5,8).Value = "Include"

$ws3.Cells.Item(6,1).Value = $fileConfig
$ws3.Cells.Item(6,2).Value = "Not to be localized"
$ws3.Cells.Item(6,4).Value = $epoch
$ws3.Cells.Item(6,7).Value = $epoch
$ws3.Cells.Item(6,8).Value = "Ignored"

$ws3.Cells.Item(1,1).Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Cells.Item(2,1), (E2eUrl $file98dc), "", "", $file98dc) | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(2,3), (DeUrl $xlf98dcDe), "", "", $xlf98dcDe) | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(3,1), (E2eUrl $fileF0b8), "", "", $fileF0b8) | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(3,3), (DeUrl $xlfF0b8De), "", "", $xlfF0b8De) | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(4,1), (E2eUrl $file6172), "", "", $file6172) | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(4,3), (DeUrl $xlf6172De), "", "", $xlf6172De) | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(5,1), (E2eUrl $file9777), "", "", $file9777) | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(5,3), (DeUrl $xlf9777De), "", "", $xlf9777De) | Out-Null
$ws3.Hyperlinks.Add($ws3.Cells.Item(6,1), (ConfigUrl), "", "", $fileConfig) | Out-Null
